# Regenerate merged AHB files
# 1) Rename header row labels: "..._old" -> "..._FV2404", "..._new" -> "..._FV2410"
# 2) Turn A1:U64 into an Excel Table ("Table1")
# 3) Freeze the header row (top row) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header cells --------------------------------------------------
# Columns A..J: "*_old" -> "*_FV2404"; columns L..U: "*_new" -> "*_FV2410"
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# Columns A..J hold the "*_old" -> "*_FV2404" headers
for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $fv2404Headers[$i]
}

# Column K holds "diff" (unchanged)

# Columns L..U hold the "*_new" -> "*_FV2410" headers
for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $fv2410Headers[$i]
}

# --- 2) Create the Excel table over A1:U64 -----------------------------------
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3) Freeze the header row -------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
